$d = $word.ActiveDocument

$d.Content.Find.Execute("23×84=1932", $true, $false, $false, $false, $false, $true, 1, $false, "49×98=4802", 2) | Out-Null
$d.Content.Find.Execute("62×82=5084", $true, $false, $false, $false, $false, $true, 1, $false, "55×93=5115", 2) | Out-Null
$d.Content.Find.Execute("85×13=1105", $true, $false, $false, $false, $false, $true, 1, $false, "78×73=5694", 2) | Out-Null
$d.Content.Find.Execute("71×88=6248", $true, $false, $false, $false, $false, $true, 1, $false, "40×67=2680", 2) | Out-Null
$d.Content.Find.Execute("47×81=3807", $true, $false, $false, $false, $false, $true, 1, $false, "44×42=1848", 2) | Out-Null
$d.Content.Find.Execute("35×56=1960", $true, $false, $false, $false, $false, $true, 1, $false, "74×90=6660", 2) | Out-Null
$d.Content.Find.Execute("25×61=1525", $true, $false, $false, $false, $false, $true, 1, $false, "41×45=1845", 2) | Out-Null
$d.Content.Find.Execute("52×44=2288", $true, $false, $false, $false, $false, $true, 1, $false, "56×29=1624", 2) | Out-Null
$d.Content.Find.Execute("19×19=361", $true, $false, $false, $false, $false, $true, 1, $false, "73×60=4380", 2) | Out-Null
$d.Content.Find.Execute("83×74=6142", $true, $false, $false, $false, $false, $true, 1, $false, "26×72=1872", 2) | Out-Null
$d.Content.Find.Execute("81×49=3969", $true, $false, $false, $false, $false, $true, 1, $false, "30×90=2700", 2) | Out-Null
$d.Content.Find.Execute("34×66=2244", $true, $false, $false, $false, $false, $true, 1, $false, "49×93=4557", 2) | Out-Null
$d.Content.Find.Execute("97×58=5626", $true, $false, $false, $false, $false, $true, 1, $false, "77×30=2310", 2) | Out-Null
$d.Content.Find.Execute("33×62=2046", $true, $false, $false, $false, $false, $true, 1, $false, "98×66=6468", 2) | Out-Null
$d.Content.Find.Execute("20×77=1540", $true, $false, $false, $false, $false, $true, 1, $false, "17×38=646", 2) | Out-Null
$d.Content.Find.Execute("73×41=2993", $true, $false, $false, $false, $false, $true, 1, $false, "90×83=7470", 2) | Out-Null
$d.Content.Find.Execute("63×64=4032", $true, $false, $false, $false, $false, $true, 1, $false, "52×20=1040", 2) | Out-Null
$d.Content.Find.Execute("92×20=1840", $true, $false, $false, $false, $false, $true, 1, $false, "27×30=810", 2) | Out-Null
$d.Content.Find.Execute("15×55=825", $true, $false, $false, $false, $false, $true, 1, $false, "27×57=1539", 2) | Out-Null
$d.Content.Find.Execute("73×51=3723", $true, $false, $false, $false, $false, $true, 1, $false, "90×30=2700", 2) | Out-Null
$d.Content.Find.Execute("43×47=2021", $true, $false, $false, $false, $false, $true, 1, $false, "57×12=684", 2) | Out-Null
$d.Content.Find.Execute("69×16=1104", $true, $false, $false, $false, $false, $true, 1, $false, "86×57=4902", 2) | Out-Null
$d.Content.Find.Execute("46×26=1196", $true, $false, $false, $false, $false, $true, 1, $false, "69×26=1794", 2) | Out-Null
$d.Content.Find.Execute("26×45=1170", $true, $false, $false, $false, $false, $true, 1, $false, "17×43=731", 2) | Out-Null
$d.Content.Find.Execute("61×86=5246", $true, $false, $false, $false, $false, $true, 1, $false, "36×23=828", 2) | Out-Null
